# Fruta / hortaliza, semanal
# Insert two new price-record rows for "Superior Seedless" grapes
# (EE.UU. origin) at row 456, pushing the existing rows 456:537 down
# to 458:539.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 456 (shifts everything at/after 456 down by 2)
$ws.Rows("456:457").Insert()

# --- New row 456 ---
$ws.Cells.Item(456,1).Value  = 10
$ws.Cells.Item(456,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(456,3).Value  = "La Araucanía"
$ws.Cells.Item(456,4).Value  = 44504
$ws.Cells.Item(456,5).Value  = 9
$ws.Cells.Item(456,6).Value  = "Fruta"
$ws.Cells.Item(456,7).Value  = 100109
$ws.Cells.Item(456,8).Value  = "Uva"
$ws.Cells.Item(456,9).Value  = 100109001
$ws.Cells.Item(456,10).Value = "Uva"
$ws.Cells.Item(456,11).Value = "Superior Seedless"
$ws.Cells.Item(456,12).Value = "Especial"
$ws.Cells.Item(456,13).Value = 200
$ws.Cells.Item(456,14).Value = 38000
$ws.Cells.Item(456,15).Value = 38000
$ws.Cells.Item(456,16).Value = 38000
$ws.Cells.Item(456,17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(456,18).Value = "EE.UU."
$ws.Cells.Item(456,19).Value = 4750
$ws.Cells.Item(456,20).Value = 8

# --- New row 457 ---
$ws.Cells.Item(457,1).Value  = 10
$ws.Cells.Item(457,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(457,3).Value  = "La Araucanía"
$ws.Cells.Item(457,4).Value  = 44504
$ws.Cells.Item(457,5).Value  = 9
$ws.Cells.Item(457,6).Value  = "Fruta"
$ws.Cells.Item(457,7).Value  = 100109
$ws.Cells.Item(457,8).Value  = "Uva"
$ws.Cells.Item(457,9).Value  = 100109001
$ws.Cells.Item(457,10).Value = "Uva"
$ws.Cells.Item(457,11).Value = "Superior Seedless"
$ws.Cells.Item(457,12).Value = "Primera"
$ws.Cells.Item(457,13).Value = 1400
$ws.Cells.Item(457,14).Value = 33000
$ws.Cells.Item(457,15).Value = 34000
$ws.Cells.Item(457,16).Value = 33429
$ws.Cells.Item(457,17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(457,18).Value = "EE.UU."
$ws.Cells.Item(457,19).Value = 4179
$ws.Cells.Item(457,20).Value = 8
